{"js": "// The template used a Word field (fldChar begin/instrText.../fldChar end)\n// to hold the M2Doc expression  { m:'doc.html'.fromHTMLURI() } .\n// The parser was switched to TokenIteratorFieldRewriterSplit, which expects\n// that same expression written as plain literal text (with literal curly\n// braces) instead of as a Word field. Replace that paragraph's field-code\n// runs with plain-text runs carrying the equivalent characters, keeping the\n// bookmark (\"_GoBack\") that sits in the middle of the expression untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Grab the OOXML of every paragraph so we can find the one holding the\n// field (it is the only paragraph containing a <w:fldChar> run).\nconst ooxmlResults = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  ooxmlResults.push(paragraphs.items[i].getOoxml());\n}\nawait context.sync();\n\nlet fieldParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (ooxmlResults[i].value.indexOf(\"fldChar\") !== -1) {\n    fieldParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!fieldParagraph) {\n  throw new Error(\"Could not find the paragraph containing the Word field.\");\n}\n\n// Build the replacement paragraph: the same run sequence, but as literal\n// text runs, wrapped with \"{\" / \"}\" in place of the field delimiters, and\n// with the original bookmark preserved in its original position.\nconst newParagraphOoxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:r><w:t>{</w:t></w:r>' +\n  '<w:r><w:t>m</w:t></w:r>' +\n  '<w:r><w:t>:</w:t></w:r>' +\n  \"<w:r><w:t>'</w:t></w:r>\" +\n  '<w:r><w:t>doc.html</w:t></w:r>' +\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n  '<w:bookmarkEnd w:id=\"0\"/>' +\n  \"<w:r><w:t>'.fromHTMLURI()</w:t></w:r>\" +\n  '<w:r><w:t xml:space=\"preserve\">}</w:t></w:r>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\nfieldParagraph.insertOoxml(newParagraphOoxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The template used a Word field (fldChar begin/instrText.../fldChar end)\n# to hold the M2Doc expression  { m:'doc.html'.fromHTMLURI() } .\n# The parser was switched to TokenIteratorFieldRewriterSplit, which expects\n# that same expression written as plain literal text (with literal curly\n# braces) instead of as a Word field. Replace that paragraph's field-code\n# runs with plain-text runs carrying the equivalent characters, keeping the\n# bookmark (\"_GoBack\") that sits in the middle of the expression untouched.\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph that holds the field (the only one with a fldChar).\n$targetIndex = -1\nif ($d.Fields.Count -gt 0) {\n    $f = $d.Fields.Item(1)\n    $fieldRange = $f.Code\n    $paras = $d.Paragraphs\n    for ($i = 1; $i -le $paras.Count; $i++) {\n        $p = $paras.Item($i)\n        if ($p.Range.Start -le $fieldRange.Start -and $p.Range.End -ge $fieldRange.End) {\n            $targetIndex = $i\n            break\n        }\n    }\n}\n\nif ($targetIndex -eq -1) {\n    throw \"Could not find the paragraph containing the Word field.\"\n}\n\n$targetParagraph = $d.Paragraphs.Item($targetIndex)\n$targetRange = $targetParagraph.Range\n\n# Replace the whole paragraph's XML: the same run sequence, but as literal\n# text runs, wrapped with \"{\" / \"}\" in place of the field delimiters, and\n# with the original bookmark preserved in its original position.\n$newParagraphXml = \"<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>\" +\n    \"<w:r><w:t>{</w:t></w:r>\" +\n    \"<w:r><w:t>m</w:t></w:r>\" +\n    \"<w:r><w:t>:</w:t></w:r>\" +\n    \"<w:r><w:t>'</w:t></w:r>\" +\n    \"<w:r><w:t>doc.html</w:t></w:r>\" +\n    \"<w:bookmarkStart w:id='0' w:name='_GoBack'/>\" +\n    \"<w:bookmarkEnd w:id='0'/>\" +\n    \"<w:r><w:t>'.fromHTMLURI()</w:t></w:r>\" +\n    \"<w:r><w:t xml:space='preserve'>}</w:t></w:r>\" +\n    \"</w:p>\"\n\n$targetRange.InsertXML($newParagraphXml)\n"}
